$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("Z1").Value = "hello"
$ws.Range("Z1").Style = $wb.Styles.Item(6)
Write-Host ("Z1 style name=" + $ws.Range("Z1").Style.Name)
